$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Against Rand / ml_stripped) - fill in B4:S4
$ws.Range("B4").Value = 98
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 90
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 95
$ws.Range("G4").Value = 25
$ws.Range("H4").Value = 101
$ws.Range("I4").Value = 19
$ws.Range("J4").Value = 96
$ws.Range("K4").Value = 24
$ws.Range("L4").Value = 91
$ws.Range("M4").Value = 29
$ws.Range("N4").Value = 88
$ws.Range("O4").Value = 32
$ws.Range("P4").Value = 93
$ws.Range("Q4").Value = 27
$ws.Range("R4").Value = 91
$ws.Range("S4").Value = 29

# Row 9 (Against Rand / ml_combined) - fill in B9:S9
$ws.Range("B9").Value = 107
$ws.Range("C9").Value = 13
$ws.Range("D9").Value = 105
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 107
$ws.Range("G9").Value = 13
$ws.Range("H9").Value = 106
$ws.Range("I9").Value = 14
$ws.Range("J9").Value = 107
$ws.Range("K9").Value = 13
$ws.Range("L9").Value = 103
$ws.Range("M9").Value = 17
$ws.Range("N9").Value = 99
$ws.Range("O9").Value = 21
$ws.Range("P9").Value = 107
$ws.Range("Q9").Value = 13
$ws.Range("R9").Value = 108
$ws.Range("S9").Value = 12

# Row 26 (Against Rdeep / ml_stripped) - fill in remaining N26:S26
$ws.Range("N26").Value = 45
$ws.Range("O26").Value = 75
$ws.Range("P26").Value = 44
$ws.Range("Q26").Value = 76
$ws.Range("R26").Value = 33
$ws.Range("S26").Value = 87

# Row 31 (Against Rdeep / ml_combined) - fill in remaining N31:S31
$ws.Range("N31").Value = 26
$ws.Range("O31").Value = 94
$ws.Range("P31").Value = 28
$ws.Range("Q31").Value = 92
$ws.Range("R31").Value = 61
$ws.Range("S31").Value = 59

# Update the active cell selection on the frozen (topRight) pane to S9
$ws.Range("S9").Select()
